$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log the missing "Home" hours for the week of 43143 (row 5) in column H.
$ws.Range("H5").Value = 7.25

# Reflect the cell the author was last working in (selection moved to M10).
$ws.Range("M10").Select()
